$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 22 (row 31): DC vs RCB - fill in the prediction scores
$ws.Range("E31").Value = 40
$ws.Range("H31").Value = 0
$ws.Range("K31").Value = 60
$ws.Range("N31").Value = 80
$ws.Range("Q31").Value = 20
$ws.Range("T31").Value = 100

$wb.Save()
